# Admin page now requires admin auth
#
# Fill in the "Comments / Code Files" (column E) for the implemented
# checklist items with the relevant jsp/html file names that satisfy each
# requirement, widen column E to fit the new text, and update the saved
# view/selection state to where the edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> comment/code-file text for column E
$values = @{
    12 = "listprod.jsp"
    13 = "listprod.jsp"
    14 = "listprod.jsp"
    15 = "product.jsp"
    16 = "HTML on all pages"
    19 = "Very pretty"
    23 = "listprod.jsp"
    24 = "showcart.jsp"
    25 = "showcart.jsp"
    26 = "showcart.jsp"
    31 = "checkout.jsp"
    38 = "product.jsp"
    39 = "product.jsp"
    44 = "register.html"
    45 = "processRegistration.jsp"
    46 = "customer.jsp"
    47 = "login.jsp/logout.jsp"
    52 = "product.jsp/processReview.jsp"
    53 = "product.jsp"
    54 = "processReview.jsp"
    57 = "product.jsp"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}

# Widen column E (5) to comfortably fit the new comments/code-file text.
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# The longer entries wrap to two lines at this column width, so those rows
# grow taller (matches Excel's automatic row-height adjustment for wrapped
# text).
$ws.Rows.Item(45).RowHeight = 29.5
$ws.Rows.Item(52).RowHeight = 29.5

# Update the saved scroll position / selection to reflect where the user was
# last working when the edits were made.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws.Range("E57").Select()
